$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column L values for row 4 (years header)
$ws.Range("L4").Value = 2022

# Row 5 data
$ws.Range("L5").Value = 8800.6

# Row 6 is an empty spacer row - leave L6 blank but it will get touched by style via Range below

# Row 7: formula
$ws.Range("L7").Formula = "=L5-L8"

# Row 8
$ws.Range("L8").Value = 258.39999999999998

# Row 10 - 17 data
$ws.Range("L10").Value = 683.8
$ws.Range("L11").Value = 1101.8
$ws.Range("L12").Value = 714.9
$ws.Range("L13").Value = 757.9
$ws.Range("L14").Value = 1383.3
$ws.Range("L15").Value = 1023.7
$ws.Range("L16").Value = 2929.3
$ws.Range("L17").Value = 148.9

# Row 18
$ws.Range("L18").Value = 57

# Copy styles from column K to column L for rows 4-18
$ws.Range("K4:K18").Copy()
$ws.Range("L4:L18").PasteSpecial(-4122)  # xlPasteFormats

# Update selection to reflect new active cell
$ws.Range("M4").Select()
